# Update cryptos list figures (prices / 1h volume %) and reorder a few
# coin rows, matching the latest scrape pulled in by the GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '''68.236.03'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +8.15%  '
# Row 3
$ws.Range('D3').Value = '''3.630.27'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +4.11%  '
# Row 4
$ws.Range('D4').Value = '''0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.05%  '
# Row 5
$ws.Range('D5').Value = '''418.84'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.13%  '
# Row 6
$ws.Range('D6').Value = '''133.17'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.17%  '
# Row 7
$ws.Range('D7').Value = '''0.650'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.79%  '
# Row 8
$ws.Range('D8').Value = '''3.621.28'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.05%  '
# Row 9
$ws.Range('E9').Value = '  -0.06%  '
# Row 10
$ws.Range('D10').Value = '''0.771'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.96%  '
# Row 11
$ws.Range('D11').Value = '''0.183'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +18.75%  '
# Row 12
$ws.Range('D12').Value = '''0.0000358'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +58.60%  '
# Row 13
$ws.Range('D13').Value = '''43.04'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.11%  '
# Row 14
$ws.Range('D14').Value = '''9.91'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.28%  '
# Row 15
$ws.Range('D15').Value = '''4.204.41'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.14%  '
# Row 16
$ws.Range('E16').Value = '  -0.27%  '
# Row 17
$ws.Range('D17').Value = '''20.45'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.26%  '
# Row 18
$ws.Range('D18').Value = '''3.627.94'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.66%  '
# Row 19
$ws.Range('E19').Value = '  +5.33%  '
# Row 20
$ws.Range('D20').Value = '''68.010.69'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +7.95%  '
# Row 21
$ws.Range('D21').Value = '''12.36'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.09%  '
# Row 22
$ws.Range('D22').Value = '''469.10'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.05%  '
# Row 23
$ws.Range('D23').Value = '''88.83'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.97%  '
# Row 24
$ws.Range('D24').Value = '''3.14'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -5.46%  '
# Row 25
$ws.Range('D25').Value = '''13.39'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.72%  '
# Row 26
$ws.Range('D26').Value = '''3.35'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.18%  '
# Row 27
$ws.Range('D27').Value = '''36.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +7.54%  '
# Row 28
$ws.Range('D28').Value = '''10.05'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.17%  '
# Row 29
$ws.Range('D29').Value = '''4.88'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.84%  '
# Row 30
$ws.Range('B30').Value = 'Cosmos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D30').Value = '''12.51'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.90%  '
# Row 31
$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D31').Value = '''2.78'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.96%  '
# Row 32
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').Value = '''0.118'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.06%  '
# Row 33
$ws.Range('B33').Value = 'RenderToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D33').Value = '''7.37'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.34%  '
# Row 34
$ws.Range('D34').Value = '''0.162'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.27%  '
# Row 35
$ws.Range('D35').Value = '''40.74'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.13%  '
# Row 36
$ws.Range('E36').Value = '  +0.04%  '
# Row 37
$ws.Range('D37').Value = '''56.83'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.29%  '
# Row 38
$ws.Range('D38').Value = '''0.0495'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.12%  '
# Row 39
$ws.Range('B39').Value = 'PEPE'
$ws.Range('C39').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D39').Value = '''0.0₃0704'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +18.96%  '
# Row 40
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D40').Value = '''0.146'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +7.70%  '
# Row 41
$ws.Range('D41').Value = '''0.996'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.20%  '
# Row 42
$ws.Range('D42').Value = '''3.04'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.18%  '
# Row 43
$ws.Range('D43').Value = '''148.51'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.50%  '
# Row 44
$ws.Range('E44').Value = '  -3.82%  '
# Row 45
$ws.Range('E45').Value = '  -1.16%  '
# Row 46
$ws.Range('D46').Value = '''4.34'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.02%  '
# Row 47
$ws.Range('E47').Value = '  -3.57%  '
# Row 48
$ws.Range('E48').Value = '  -3.18%  '
# Row 49
$ws.Range('D49').Value = '''2.35'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.28%  '
# Row 50
$ws.Range('D50').Value = '''2.68'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +15.20%  '
# Row 51
$ws.Range('D51').Value = '''15.68'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.36%  '
